$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential notice date from 2021-03-17 to 2021-03-18
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-9
$ws.Range("D2").Value = 0.1497630312786706
$ws.Range("E2").Value = -0.004409171075837715

$ws.Range("D3").Value = 0.1505194102245225
$ws.Range("E3").Value = -0.005025125628140503

$ws.Range("D4").Value = 0.3006606309761191
$ws.Range("E4").Value = -0.003527336860670194

$ws.Range("D5").Value = 0.1494548768933235
$ws.Range("E5").Value = -0.004686035613870709

$ws.Range("D6").Value = 0.07974395171981162
$ws.Range("E6").Value = -0.00101626016260159

$ws.Range("D7").Value = 0.1199480939950916
$ws.Range("E7").Value = -0.0009842519685039353

$ws.Range("D8").Value = 0.04991000491246115
$ws.Range("E8").Value = -0.005512679162072875

$ws.Range("E9").Value = -0.003651829566612697
